$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255308747291565
$ws.Range("B1").Value = -1
$ws.Range("D1").Value = 0.8778237700462341
$ws.Range("E1").Value = 0.9975269436836243
